# Scheduled data refresh: update cryptocurrency price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 47/48 swap coin identity (Coin name + Link) as rankings changed
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

# Update Price (column D) and Volume(1h) (column E) for every coin row
# Price values are plain text in the sheet; a leading apostrophe forces
# Excel to store them as text instead of coercing to a number/date.
$ws.Range("D2").Value = '''69.003.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.44%  '
$ws.Range("D3").Value = '''3.710.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '''615.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.14%  '
$ws.Range("D6").Value = '''191.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.04%  '
$ws.Range("D7").Value = '''0.635'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.714'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").Value = '''0.160'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.25%  '
$ws.Range("D11").Value = '''56.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.98%  '
$ws.Range("D12").Value = '''0.0000289'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.63%  '
$ws.Range("D13").Value = '''10.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '''4.296.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '''3.708.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '''0.127'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '''19.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '''1.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '''12.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '''68.762.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = '''409.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("D22").Value = '''4.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.62%  '
$ws.Range("D23").Value = '''89.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").Value = '''3.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = '''12.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").Value = '''10.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("D27").Value = '''6.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("D28").Value = '''3.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("D29").Value = '''9.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("D30").Value = '''32.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("D31").Value = '''7.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.04%  '
$ws.Range("D32").Value = '''12.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = '''0.121'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("D34").Value = '''624.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.87%  '
$ws.Range("D35").Value = '''44.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.35%  '
$ws.Range("D36").Value = '''65.66'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = '''0.411'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.25%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = '''0.0₃0807'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -11.06%  '
$ws.Range("D40").Value = '''1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").Value = '''0.140'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.27%  '
$ws.Range("D42").Value = '''3.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").Value = '''0.0441'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("D44").Value = '''2.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").Value = '''0.140'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.81%  '
$ws.Range("D46").Value = '''2.867.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.35%  '
$ws.Range("D47").Value = '''9.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.42%  '
$ws.Range("D48").Value = '''2.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D49").Value = '''3.14'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").Value = '''141.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '''2.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.66%  '
